# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the newly scraped counts.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 174
$ws1.Range("F3").Value = 660
$ws1.Range("F5").Value = 221
$ws1.Range("F6").Value = 1605
$ws1.Range("F8").Value = 3163
$ws1.Range("F9").Value = 454
$ws1.Range("F10").Value = 741

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 174
$ws4.Range("F3").Value = 660
$ws4.Range("F6").Value = 221
$ws4.Range("F7").Value = 1605
$ws4.Range("F9").Value = 3163
$ws4.Range("F10").Value = 454
$ws4.Range("F11").Value = 741
